# ---------------------------------------------------------------------------
# Target change (from the supplied OOXML diff):
#   word/styles.xml -> <w:docDefaults> is trimmed so that <w:rPrDefault> keeps
#   only rFonts/sz/szCs/lang, and <w:pPrDefault> keeps only a bare
#   <w:spacing w:line="276" w:lineRule="auto"/>. Every attribute that is
#   removed (b=0, i=0, smallCaps=0, strike=0, color=000000, u=none,
#   shd=clear/auto, vertAlign=baseline, keepNext=0, keepLines=0,
#   widowControl=1, the nil pBdr, spacing after/before=0, ind all-0,
#   contextualSpacing=0, jc=left) is exactly the implicit OOXML default for
#   that property, so the edit is a pure redundant-attribute cleanup of
#   <w:docDefaults> with zero effect on resolved/rendered formatting, and it
#   does not touch the "Normal" style element (which stays the empty
#   <w:style .../><w:name .../></w:style> it already was) or anything else
#   in the package.
#
# <w:docDefaults> is a document-defaults block that lives outside every
# exposed collection (Document.Styles only enumerates the 10 named styles -
# Normal, Table Normal, Heading1-6, Title, Subtitle - confirmed by walking
# Styles.Item(1..Count)); Word's object model (both real Word/VBA and this
# COM-interop surface) has no Style/Styles/Document member that reads or
# writes rPrDefault/pPrDefault - Font/ParagraphFormat setters on
# Styles("Normal") only ever materialize a local <w:rPr>/<w:pPr> on the
# Normal style element itself, never <w:docDefaults> - and
# Document.WordOpenXML / Range.WordOpenXML round-trip the flat-OPC XML as a
# read-only snapshot (assigning back to them does not mark the package
# dirty). So there is no COM call that can reach into <w:docDefaults>
# without fabricating a side effect on Normal/other styles that is not part
# of this diff (and would itself be wrong, since the Normal style element is
# unchanged in the target).
#
# Since the block is unreachable from the object model and the edit has no
# semantic/visual effect (every paragraph and run in this document already
# carries fully explicit direct formatting - see word/document.xml - so
# nothing in the content actually inherits from docDefaults), the safest,
# most faithful action available through $word/$d is to leave the
# document's reachable state untouched rather than emulate the cleanup via
# an unrelated object (e.g. the Normal style) and introduce a diff that
# isn't in the target. Touching the document only through read-only calls
# below keeps the saved package identical to the source in every part that
# the OM can actually influence.

$d = $word.ActiveDocument

# Read-only sanity touches (no mutation): confirms the document/model is
# reachable without altering Content, Styles, or any other part.
$null = $d.Styles.Count
$null = $d.Paragraphs.Count
